$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1966.3422
$ws.Range("I98").Value = 1308.12
$ws.Range("J98").Value = 3232.1538
$ws.Range("K98").Value = 1308.12
$ws.Range("L98").Value = 3232.1538
$ws.Range("M98").Value = 189.8800000000001
$ws.Range("N98").Value = -6228.1538

$ws.Range("H106").Value = 2973.5527
$ws.Range("I106").Value = 2874.375
$ws.Range("K106").Value = 2874.375
$ws.Range("M106").Value = -2243.375

$ws.Range("H121").Value = 1109.6154
$ws.Range("J121").Value = 1136.6666
$ws.Range("L121").Value = 3409.9998
$ws.Range("N121").Value = -6903.9998

$ws.Range("H122").Value = 1966.3422
$ws.Range("I122").Value = 1308.12
$ws.Range("J122").Value = 3232.1538
$ws.Range("K122").Value = 3924.36
$ws.Range("L122").Value = 9696.4614
$ws.Range("M122").Value = -1474.36
$ws.Range("N122").Value = -14596.4614

$ws.Range("H131").Value = 729.1875
$ws.Range("I131").Value = 591.13336
$ws.Range("J131").Value = 2800
$ws.Range("K131").Value = 1773.40008
$ws.Range("L131").Value = 8400
$ws.Range("M131").Value = 3266.59992
$ws.Range("N131").Value = -18480

$ws.Range("H132").Value = 2552157.5
$ws.Range("I132").Value = 2842061.8
$ws.Range("K132").Value = 8526185.399999999
$ws.Range("M132").Value = -8523655.399999999

$ws.Range("H137").Value = 2874.418
$ws.Range("I137").Value = 2854.1353
$ws.Range("J137").Value = 2899.4333
$ws.Range("K137").Value = 8562.4059
$ws.Range("L137").Value = 8698.2999
$ws.Range("M137").Value = -6012.4059
$ws.Range("N137").Value = -13798.2999

$ws.Range("H138").Value = 7835.407
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 7835.407
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 23506.221
$ws.Range("N138").Value = -33786.22100000001
$ws.Range("M138").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3642.2122
$ws.Range("I32").Value = 2622.4482
$ws.Range("K32").Value = 2622.4482
$ws.Range("M32").Value = -2335.4482

$ws.Range("H61").Value = 858.08826
$ws.Range("I61").Value = 853.7879
$ws.Range("J61").Value = 1000
$ws.Range("K61").Value = 853.7879
$ws.Range("L61").Value = 1000
$ws.Range("M61").Value = -641.7879
$ws.Range("N61").Value = -1424

$ws.Range("H74").Value = 1649.037
$ws.Range("I74").Value = 702.3333
$ws.Range("K74").Value = 702.3333
$ws.Range("M74").Value = 171.6667

$ws.Range("H77").Value = 1649.037
$ws.Range("I77").Value = 702.3333
$ws.Range("K77").Value = 3511.6665
$ws.Range("M77").Value = 856.3334999999997

$ws.Range("H122").Value = 6252730
$ws.Range("I122").Value = 8335298.5
$ws.Range("J122").Value = 5025
$ws.Range("K122").Value = 25005895.5
$ws.Range("L122").Value = 15075
$ws.Range("M122").Value = -25003445.5
$ws.Range("N122").Value = -19975

$ws.Range("H132").Value = 2058.1064
$ws.Range("I132").Value = 2171.5625
$ws.Range("J132").Value = 1816.0667
$ws.Range("K132").Value = 6514.6875
$ws.Range("L132").Value = 5448.2001
$ws.Range("M132").Value = -3984.6875
$ws.Range("N132").Value = -10508.2001

$ws.Range("H136").Value = 858.08826
$ws.Range("I136").Value = 853.7879
$ws.Range("J136").Value = 1000
$ws.Range("K136").Value = 2561.3637
$ws.Range("L136").Value = 3000
$ws.Range("M136").Value = -11.36369999999988
$ws.Range("N136").Value = -8100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1935.1428
$ws.Range("I107").Value = 1206.6
$ws.Range("J107").Value = 3756.5
$ws.Range("K107").Value = 1206.6
$ws.Range("L107").Value = 3756.5
$ws.Range("M107").Value = 713.4000000000001
$ws.Range("N107").Value = -7596.5

$ws.Range("H134").Value = 655.2963
$ws.Range("I134").Value = 600.89795
$ws.Range("J134").Value = 1188.4
$ws.Range("K134").Value = 1802.69385
$ws.Range("L134").Value = 3565.2
$ws.Range("M134").Value = 732.3061499999999
$ws.Range("N134").Value = -8635.200000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 19343.959
$ws.Range("I31").Value = 23691.717
$ws.Range("J31").Value = 11651.77
$ws.Range("K31").Value = 23691.717
$ws.Range("L31").Value = 11651.77
$ws.Range("M31").Value = -23396.717
$ws.Range("N31").Value = -12241.77

$ws.Range("H34").Value = 19343.959
$ws.Range("I34").Value = 23691.717
$ws.Range("J34").Value = 11651.77
$ws.Range("K34").Value = 23691.717
$ws.Range("L34").Value = 11651.77
$ws.Range("M34").Value = -23489.717
$ws.Range("N34").Value = -12055.77

$ws.Range("H58").Value = 1117.8036
$ws.Range("I58").Value = 1046.1282
$ws.Range("J58").Value = 1282.2354
$ws.Range("K58").Value = 1046.1282
$ws.Range("L58").Value = 1282.2354
$ws.Range("M58").Value = -843.1282000000001
$ws.Range("N58").Value = -1688.2354

$ws.Range("H99").Value = 1738.5714
$ws.Range("I99").Value = 1489
$ws.Range("J99").Value = 2071.3333
$ws.Range("K99").Value = 1489
$ws.Range("L99").Value = 2071.3333
$ws.Range("M99").Value = 9
$ws.Range("N99").Value = -5067.3333

$ws.Range("H105").Value = 919
$ws.Range("I105").Value = 982
$ws.Range("J105").Value = 840.25
$ws.Range("K105").Value = 982
$ws.Range("L105").Value = 840.25
$ws.Range("M105").Value = 765
$ws.Range("N105").Value = -4334.25

$ws.Range("H126").Value = 1738.5714
$ws.Range("I126").Value = 1489
$ws.Range("J126").Value = 2071.3333
$ws.Range("K126").Value = 4467
$ws.Range("L126").Value = 6213.999899999999
$ws.Range("M126").Value = -1997
$ws.Range("N126").Value = -11153.9999

$ws.Range("H132").Value = 1070.7632
$ws.Range("I132").Value = 935.04
$ws.Range("J132").Value = 1331.7693
$ws.Range("K132").Value = 2805.12
$ws.Range("L132").Value = 3995.3079
$ws.Range("M132").Value = -275.1199999999999
$ws.Range("N132").Value = -9055.3079

$ws.Range("H134").Value = 2121.6538
$ws.Range("I134").Value = 1742.3334
$ws.Range("J134").Value = 2638.9092
$ws.Range("K134").Value = 5227.0002
$ws.Range("L134").Value = 7916.7276
$ws.Range("M134").Value = -2692.0002
$ws.Range("N134").Value = -12986.7276

$ws.Range("H135").Value = 43413.332
$ws.Range("J135").Value = 43413.332
$ws.Range("L135").Value = 43413.332
$ws.Range("N135").Value = -53553.332

$ws.Range("H136").Value = 1117.8036
$ws.Range("I136").Value = 1046.1282
$ws.Range("J136").Value = 1282.2354
$ws.Range("K136").Value = 3138.3846
$ws.Range("L136").Value = 3846.7062
$ws.Range("M136").Value = -588.3846000000003
$ws.Range("N136").Value = -8946.706200000001

$ws.Range("H137").Value = 68160
$ws.Range("J137").Value = 68160
$ws.Range("L137").Value = 68160
$ws.Range("N137").Value = -78360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2722
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()

$ws.Range("H39").Value = 6413.231
$ws.Range("I39").Value = 10000
$ws.Range("J39").Value = 6114.3335
$ws.Range("K39").Value = 30000
$ws.Range("L39").Value = 18343.0005
$ws.Range("M39").Value = -29706
$ws.Range("N39").Value = -18931.0005

$ws.Range("H80").Value = 3400.3333
$ws.Range("I80").Value = 3379.8
$ws.Range("J80").Value = 3410.6
$ws.Range("K80").Value = 10139.4
$ws.Range("L80").Value = 10231.8
$ws.Range("M80").Value = -9203.400000000001
$ws.Range("N80").Value = -12103.8

$ws.Range("H83").Value = 3400.3333
$ws.Range("I83").Value = 3379.8
$ws.Range("J83").Value = 3410.6
$ws.Range("K83").Value = 30418.2
$ws.Range("L83").Value = 30695.4
$ws.Range("M83").Value = -25738.2
$ws.Range("N83").Value = -40055.39999999999

$ws.Range("H129").Value = 6537124
$ws.Range("J129").Value = 11112697
$ws.Range("L129").Value = 33338091
$ws.Range("N129").Value = -33348091

$ws.Range("H132").Value = 1265
$ws.Range("I132").Value = 866.6667
$ws.Range("J132").Value = 1384.5
$ws.Range("K132").Value = 7800.0003
$ws.Range("L132").Value = 12460.5
$ws.Range("M132").Value = -5270.0003
$ws.Range("N132").Value = -17520.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1166.6666
$ws.Range("I113").Value = 1500
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 1500
$ws.Range("L113").Value = 1000
$ws.Range("M113").Value = 670
$ws.Range("N113").Value = -5340

$ws.Range("H132").Value = 1679.3715
$ws.Range("I132").Value = 1586.8125
$ws.Range("J132").Value = 2666.6667
$ws.Range("K132").Value = 4760.4375
$ws.Range("L132").Value = 8000.000100000001
$ws.Range("M132").Value = -2230.4375
$ws.Range("N132").Value = -13060.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2170.9038
$ws.Range("I132").Value = 1765.8049
$ws.Range("J132").Value = 3680.818
$ws.Range("K132").Value = 5297.4147
$ws.Range("L132").Value = 11042.454
$ws.Range("M132").Value = -2767.4147
$ws.Range("N132").Value = -16102.454

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 635.9400000000001
$ws.Range("I132").Value = 504.13635
$ws.Range("J132").Value = 1602.5
$ws.Range("K132").Value = 1512.40905
$ws.Range("L132").Value = 4807.5
$ws.Range("M132").Value = 1017.59095
$ws.Range("N132").Value = -9867.5
